# Insert a new weekly record row for "Feria Lagunitas de Puerto Montt - Mango"
# at row 267, pushing the existing rows 267-308 down to 268-309.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 267 (shifts 267..308 -> 268..309)
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with the new data record
$ws.Range("A267").Value() = 4
$ws.Range("B267").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value() = "Los Lagos"
$ws.Range("D267").Value() = 44951
$ws.Range("E267").Value() = 10
$ws.Range("F267").Value() = "Fruta"
$ws.Range("G267").Value() = 100108
$ws.Range("H267").Value() = "Tropicales y subtropicales"
$ws.Range("I267").Value() = 100108002
$ws.Range("J267").Value() = "Mango"
$ws.Range("K267").Value() = "Sin especificar"
$ws.Range("L267").Value() = "Primera"
$ws.Range("M267").Value() = 120
$ws.Range("N267").Value() = 7500
$ws.Range("O267").Value() = 8000
$ws.Range("P267").Value() = 7750
$ws.Range("Q267").Value() = "$/bandeja 4 kilos"
$ws.Range("R267").Value() = "Brasil"
$ws.Range("S267").Value() = 1938
$ws.Range("T267").Value() = 4
